$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update Marking row (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update Total row (B12): 48 -> 80
$ws.Range("B12").Value = 80

# Update Correct/Total marks text (E12): "45/84" -> "80/140"
$ws.Range("E12").Value = "80/140"
